# 22/2-2018 Made MovingObject Drawable Stable
# Append the two new lab-diary entries (19/2 and 22/2) to the bottom of the
# existing log, then move the selection to mirror where the author left the
# cursor when the workbook was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 19/2-2018 entry
$ws.Range("A5").Value = "19/2-2018"
$ws.Range("B5").Value = "Image Serching and spritesheat building"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 30

# 22/2-2018 entry
$ws.Range("A6").Value = "22/2-2018"
$ws.Range("B6").Value = "MovingObject made it drawable"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 0

# The Total tid row (row 13) sums C2:C12/D2:D12 via existing formulas, so it
# recalculates to 13:30 automatically once the new hours/minutes are in place.

# Leave the cursor on E10, matching the saved selection in the workbook.
$ws.Range("E10").Select()
